$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-validated HTTPS status (column D): rows that flipped from HTTPS to UNKNOWN
$ws.Range("D6:D10").Value = "UNKNOWN"
$ws.Range("D12:D13").Value = "UNKNOWN"
$ws.Range("D15:D17").Value = "UNKNOWN"
$ws.Range("D19:D21").Value = "UNKNOWN"
$ws.Range("D24").Value = "UNKNOWN"
$ws.Range("D31").Value = "UNKNOWN"
$ws.Range("D35").Value = "UNKNOWN"
$ws.Range("D39").Value = "UNKNOWN"
$ws.Range("D66:D68").Value = "UNKNOWN"
$ws.Range("D70").Value = "UNKNOWN"
$ws.Range("D76").Value = "UNKNOWN"
$ws.Range("D106").Value = "UNKNOWN"
$ws.Range("D115").Value = "UNKNOWN"
$ws.Range("D136").Value = "UNKNOWN"
$ws.Range("D146").Value = "UNKNOWN"
$ws.Range("D148").Value = "UNKNOWN"
$ws.Range("D196").Value = "UNKNOWN"
$ws.Range("D199").Value = "UNKNOWN"
$ws.Range("D201").Value = "UNKNOWN"
$ws.Range("D203").Value = "UNKNOWN"
$ws.Range("D220:D222").Value = "UNKNOWN"
$ws.Range("D225:D227").Value = "UNKNOWN"
$ws.Range("D232:D233").Value = "UNKNOWN"
$ws.Range("D236").Value = "UNKNOWN"
$ws.Range("D238").Value = "UNKNOWN"
$ws.Range("D241").Value = "UNKNOWN"
$ws.Range("D244:D245").Value = "UNKNOWN"

# Re-validated HTTPS status (column D): rows that flipped from UNKNOWN to HTTPS
$ws.Range("D61:D62").Value = "HTTPS"
$ws.Range("D86:D89").Value = "HTTPS"
$ws.Range("D99:D104").Value = "HTTPS"
$ws.Range("D108").Value = "HTTPS"
$ws.Range("D111:D112").Value = "HTTPS"
$ws.Range("D116:D117").Value = "HTTPS"
$ws.Range("D120:D123").Value = "HTTPS"
$ws.Range("D125:D128").Value = "HTTPS"
$ws.Range("D130:D135").Value = "HTTPS"
$ws.Range("D137:D145").Value = "HTTPS"
$ws.Range("D147").Value = "HTTPS"
$ws.Range("D149:D151").Value = "HTTPS"
$ws.Range("D153").Value = "HTTPS"
$ws.Range("D156").Value = "HTTPS"
$ws.Range("D159").Value = "HTTPS"
$ws.Range("D163").Value = "HTTPS"
$ws.Range("D165:D171").Value = "HTTPS"
$ws.Range("D176").Value = "HTTPS"
$ws.Range("D183:D184").Value = "HTTPS"
$ws.Range("D186").Value = "HTTPS"
$ws.Range("D212").Value = "HTTPS"
$ws.Range("D215").Value = "HTTPS"
$ws.Range("D250").Value = "HTTPS"

